# "add preview in datasets"
# Insert two new metadata rows (dataset.preview.table / dataset.preview.line)
# right after the "dataset.status" row (i.e. before the old row 4), each
# holding a multi-line DS-query formula string, wrapped + vertically centred,
# with a taller row height to show the full text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

$formulaTable = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nlimit(start:0, length:5);"
$formulaLine = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nline(x:-1);"

# Make room for the two new rows right after row 3 (dataset.status).
$ws.Rows("4:5").Insert()

$ws.Range("A4").Value = "dataset.preview.table"
$ws.Range("B4").Value = $formulaTable

$ws.Range("A5").Value = "dataset.preview.line"
$ws.Range("B5").Value = $formulaLine

$ws.Range("A4:B5").VerticalAlignment = -4108
$ws.Range("A4:B5").WrapText = $true

$ws.Rows(4).RowHeight = 120
$ws.Rows(5).RowHeight = 120

$ws.Range("B8").Select()
